$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 114
$ws.Range("I2").Value = 117.4
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 117.4
$ws.Range("L2").Value = 80
$ws.Range("M2").Value = -4.400000000000006
$ws.Range("N2").Value = -306
# Row 39
$ws.Range("H39").Value = 385.22223
$ws.Range("I39").Value = 108.375
$ws.Range("J39").Value = 2600
$ws.Range("K39").Value = 325.125
$ws.Range("L39").Value = 7800
$ws.Range("M39").Value = -29.125
$ws.Range("N39").Value = -8392
# Row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
# Row 53
$ws.Range("H53").Value = 121.2
$ws.Range("I53").Value = 169.66667
$ws.Range("J53").Value = 48.5
$ws.Range("K53").Value = 169.66667
$ws.Range("L53").Value = 48.5
$ws.Range("M53").Value = 467.33333
$ws.Range("N53").Value = -1322.5
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
# Row 70
$ws.Range("H70").Value = 41200.2
$ws.Range("I70").Value = 41200.2
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 123600.6
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -123330.6
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 41200.2
$ws.Range("I73").Value = 41200.2
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 123600.6
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -122664.6
$ws.Range("N73").ClearContents()
# Row 137
$ws.Range("H137").Value = 12634.458
$ws.Range("I137").Value = 2523.9092
$ws.Range("J137").Value = 21189.54
$ws.Range("K137").Value = 7571.7276
$ws.Range("L137").Value = 63568.62
$ws.Range("M137").Value = -5021.7276
$ws.Range("N137").Value = -68668.62

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 39
$ws.Range("H39").Value = 13999.75
$ws.Range("I39").Value = 8666.333000000001
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 8666.333000000001
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -8146.333000000001
$ws.Range("N39").Value = -31040
# Row 45
$ws.Range("H45").Value = 6967.591
$ws.Range("I45").Value = 7682.7896
$ws.Range("J45").Value = 2438
$ws.Range("K45").Value = 7682.7896
$ws.Range("L45").Value = 2438
$ws.Range("M45").Value = -7305.7896
$ws.Range("N45").Value = -3192
# Row 92
$ws.Range("H92").Value = 34999.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 34999.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 34999.5
$ws.Range("N92").Value = -39991.5
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 122
$ws.Range("H122").Value = 2639.2632
$ws.Range("I122").Value = 2927
$ws.Range("J122").Value = 1104.6666
$ws.Range("K122").Value = 8781
$ws.Range("L122").Value = 3313.9998
$ws.Range("M122").Value = -6331
$ws.Range("N122").Value = -8213.9998
# Row 124
$ws.Range("H124").Value = 44000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 44000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 44000
$ws.Range("N124").Value = -53820
# Row 132
$ws.Range("H132").Value = 4560.227
$ws.Range("I132").Value = 4723.6113
$ws.Range("J132").Value = 3825
$ws.Range("K132").Value = 14170.8339
$ws.Range("L132").Value = 11475
$ws.Range("M132").Value = -11640.8339
$ws.Range("N132").Value = -16535

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 912588.6
$ws.Range("I86").Value = 1669329.1
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 1669329.1
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -1668206.1
$ws.Range("N86").Value = -6746
# Row 89
$ws.Range("H89").Value = 912588.6
$ws.Range("I89").Value = 1669329.1
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 8346645.5
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -8341029.5
$ws.Range("N89").Value = -33732
# Row 94
$ws.Range("H94").Value = 352
$ws.Range("I94").Value = 340.07144
$ws.Range("J94").Value = 463.33334
$ws.Range("K94").Value = 340.07144
$ws.Range("L94").Value = 463.33334
$ws.Range("M94").Value = 110.92856
$ws.Range("N94").Value = -1365.33334
# Row 134
$ws.Range("H134").Value = 4557.9375
$ws.Range("I134").Value = 4707.7
$ws.Range("J134").Value = 4308.3335
$ws.Range("K134").Value = 14123.1
$ws.Range("L134").Value = 12925.0005
$ws.Range("M134").Value = -11588.1
$ws.Range("N134").Value = -17995.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 27500
$ws.Range("I25").Value = 25000
$ws.Range("J25").Value = 30000
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = -24826
$ws.Range("N25").Value = -30348
# Row 31
$ws.Range("H31").Value = 2708.3215
$ws.Range("I31").Value = 1794.909
$ws.Range("J31").Value = 3299.353
$ws.Range("K31").Value = 1794.909
$ws.Range("L31").Value = 3299.353
$ws.Range("M31").Value = -1499.909
$ws.Range("N31").Value = -3889.353
# Row 34
$ws.Range("H34").Value = 2708.3215
$ws.Range("I34").Value = 1794.909
$ws.Range("J34").Value = 3299.353
$ws.Range("K34").Value = 1794.909
$ws.Range("L34").Value = 3299.353
$ws.Range("M34").Value = -1592.909
$ws.Range("N34").Value = -3703.353

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 1001.5
$ws.Range("I44").Value = 1001.5
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 3004.5
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -2606.5
$ws.Range("N44").ClearContents()
# Row 117
$ws.Range("H117").Value = 590.9091
$ws.Range("I117").Value = 561.6667
$ws.Range("J117").Value = 626
$ws.Range("K117").Value = 1685.0001
$ws.Range("L117").Value = 1878
$ws.Range("M117").Value = 1756.9999
$ws.Range("N117").Value = -8762
# Row 124
$ws.Range("H124").Value = 1920
$ws.Range("I124").Value = 1920
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 5760
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -850

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 9899
$ws.Range("I20").Value = 9900
$ws.Range("J20").Value = 9897
$ws.Range("K20").Value = 9900
$ws.Range("L20").Value = 9897
$ws.Range("M20").Value = -9655
$ws.Range("N20").Value = -10387
# Row 107
$ws.Range("H107").Value = 765.7857
$ws.Range("I107").Value = 736.125
$ws.Range("J107").Value = 805.3333
$ws.Range("K107").Value = 736.125
$ws.Range("L107").Value = 805.3333
$ws.Range("M107").Value = 1183.875
$ws.Range("N107").Value = -4645.3333
# Row 122
$ws.Range("H122").Value = 3585.111
$ws.Range("I122").Value = 3443.8
$ws.Range("J122").Value = 3761.75
$ws.Range("K122").Value = 10331.4
$ws.Range("L122").Value = 11285.25
$ws.Range("M122").Value = -7881.400000000001
$ws.Range("N122").Value = -16185.25
# Row 126
$ws.Range("H126").Value = 3642.3157
$ws.Range("I126").Value = 2352.7144
$ws.Range("J126").Value = 4394.5835
$ws.Range("K126").Value = 7058.1432
$ws.Range("L126").Value = 13183.7505
$ws.Range("M126").Value = -4588.1432
$ws.Range("N126").Value = -18123.7505
# Row 132
$ws.Range("H132").Value = 3421.75
$ws.Range("I132").Value = 2977.6072
$ws.Range("J132").Value = 4976.25
$ws.Range("K132").Value = 8932.821599999999
$ws.Range("L132").Value = 14928.75
$ws.Range("M132").Value = -6402.821599999999
$ws.Range("N132").Value = -19988.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2319.9285
$ws.Range("I16").Value = 2737.9
$ws.Range("J16").Value = 1275
$ws.Range("K16").Value = 2737.9
$ws.Range("L16").Value = 1275
$ws.Range("M16").Value = -2567.9
$ws.Range("N16").Value = -1615
# Row 22
$ws.Range("H22").Value = 3406.923
$ws.Range("I22").Value = 3291.5715
$ws.Range("J22").Value = 3541.5
$ws.Range("K22").Value = 3291.5715
$ws.Range("L22").Value = 3541.5
$ws.Range("M22").Value = -2996.5715
$ws.Range("N22").Value = -4131.5
# Row 27
$ws.Range("H27").Value = 3406.923
$ws.Range("I27").Value = 3291.5715
$ws.Range("J27").Value = 3541.5
$ws.Range("K27").Value = 3291.5715
$ws.Range("L27").Value = 3541.5
$ws.Range("M27").Value = -3184.5715
$ws.Range("N27").Value = -3755.5
# Row 40
$ws.Range("H40").Value = 1799.8
$ws.Range("I40").Value = 1799.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1799.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1663.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 6541.391
$ws.Range("I126").Value = 6534
$ws.Range("J126").Value = 6562.3335
$ws.Range("K126").Value = 19602
$ws.Range("L126").Value = 19687.0005
$ws.Range("M126").Value = -17132
$ws.Range("N126").Value = -24627.0005
